$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "290.40"
Set-TextValue "E2" "-3.45%"
Set-TextValue "G2" "11"
Set-TextValue "D3" "30.70"
Set-TextValue "E3" "-5.05%"
Set-TextValue "G3" "11"
Set-TextValue "D4" "4.951"
Set-TextValue "E4" "0.02%"
Set-TextValue "G4" "11"
Set-TextValue "E5" "-5.28%"
Set-TextValue "G5" "11"
Set-TextValue "D6" "1.862"
Set-TextValue "E6" "-4.78%"
Set-TextValue "G6" "11"
Set-TextValue "D7" "7.696"
Set-TextValue "E7" "-1.83%"
Set-TextValue "G7" "11"
Set-TextValue "E8" "-0.75%"
Set-TextValue "G8" "11"
Set-TextValue "D9" "0.8969"
Set-TextValue "E9" "-2.17%"
Set-TextValue "G9" "11"
Set-TextValue "D10" "0.1658"
Set-TextValue "E10" "-5.29%"
Set-TextValue "G10" "11"
Set-TextValue "D11" "0.07730"
Set-TextValue "E11" "-0.25%"
Set-TextValue "G11" "11"
Set-TextValue "D12" "0.07945"
Set-TextValue "E12" "-7.09%"
Set-TextValue "G12" "11"
Set-TextValue "D13" "0.03032"
Set-TextValue "E13" "-5.11%"
Set-TextValue "G13" "11"
Set-TextValue "D14" "0.1002"
Set-TextValue "E14" "0.18%"
Set-TextValue "G14" "11"
Set-TextValue "D15" "0.001508"
Set-TextValue "E15" "-0.75%"
Set-TextValue "G15" "11"
Set-TextValue "D16" "0.005711"
Set-TextValue "E16" "-3.86%"
Set-TextValue "G16" "11"
Set-TextValue "G17" "11"
Set-TextValue "D18" "3.464"
Set-TextValue "E18" "0.09%"
Set-TextValue "G18" "11"
Set-TextValue "D19" "2.083"
Set-TextValue "E19" "-3.28%"
Set-TextValue "G19" "11"
Set-TextValue "D20" "0.3320"
Set-TextValue "E20" "-0.89%"
Set-TextValue "G20" "11"
Set-TextValue "D21" "0.1279"
Set-TextValue "E21" "-1.67%"
Set-TextValue "G21" "11"
Set-TextValue "D22" "4.052"
Set-TextValue "E22" "-4.99%"
Set-TextValue "G22" "11"
Set-TextValue "D23" "0.2387"
Set-TextValue "E23" "19.76%"
Set-TextValue "G23" "11"
Set-TextValue "D24" "0.04495"
Set-TextValue "E24" "-0.42%"
Set-TextValue "G24" "11"
Set-TextValue "D25" "0.001215"
Set-TextValue "E25" "-0.55%"
Set-TextValue "G25" "11"
Set-TextValue "D26" "0.004641"
Set-TextValue "E26" "5.87%"
Set-TextValue "G26" "11"
Set-TextValue "D27" "0.0001252"
Set-TextValue "E27" "0.03%"
Set-TextValue "G27" "11"
Set-TextValue "G28" "11"
Set-TextValue "G29" "11"
Set-TextValue "G30" "11"
Set-TextValue "G31" "11"
Set-TextValue "G32" "11"
Set-TextValue "G33" "11"
Set-TextValue "G34" "11"
Set-TextValue "G35" "11"
Set-TextValue "G36" "11"
Set-TextValue "G37" "11"
Set-TextValue "G38" "11"
Set-TextValue "D39" "0.01583"
Set-TextValue "E39" "-6.67%"
Set-TextValue "G39" "11"
Set-TextValue "D40" "0.04385"
Set-TextValue "E40" "-6.28%"
Set-TextValue "G40" "11"
Set-TextValue "D41" "0.007332"
Set-TextValue "E41" "-1.80%"
Set-TextValue "G41" "11"
Set-TextValue "D42" "0.009933"
Set-TextValue "G42" "11"
Set-TextValue "D43" "0.1305"
Set-TextValue "E43" "-3.07%"
Set-TextValue "G43" "11"
Set-TextValue "D44" "0.002073"
Set-TextValue "E44" "-11.13%"
Set-TextValue "G44" "11"
Set-TextValue "D45" "0.009405"
Set-TextValue "E45" "-10.43%"
Set-TextValue "G45" "11"
Set-TextValue "D46" "0.00005948"
Set-TextValue "E46" "-4.68%"
Set-TextValue "G46" "11"
Set-TextValue "E47" "-0.16%"
Set-TextValue "G47" "11"
Set-TextValue "D48" "2.254"
Set-TextValue "E48" "174.68%"
Set-TextValue "G48" "11"
Set-TextValue "G49" "11"
Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "-0.16%"
Set-TextValue "G50" "11"
Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "-0.16%"
Set-TextValue "G51" "11"
